# Update the "2024" worksheet: a new payment log entry was recorded for
# "amazeloan" (September) at 2024-09-01 10:12:03. The log keeps its newest
# entry at row 24, so recording it pushes the existing amazeloan / hdfc log
# rows (and the trailing "Broadband" marker row) down by one row each.
#
# Inserting a whole row (rather than overwriting cell-by-cell) is what
# actually shifts all the existing row data/formatting down correctly and
# keeps it intact, and leaves a clean blank row 24 for the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows("24:24").Insert()

# --- Record the new September (amazeloan) entry in the freshly opened row ---
$ws.Range("R24").Value = "amazeloan"
$ws.Range("S24").Value = "2024-09-01 10:12:03"
